# Update "want to go" counters (column F) on the "展览" and "全部类型" sheets
# to reflect the regenerated data output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 11715
$ws1.Range("F7").Value  = 11664
$ws1.Range("F9").Value  = 1166
$ws1.Range("F12").Value = 5773
$ws1.Range("F15").Value = 184
$ws1.Range("F16").Value = 18

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 11715
$ws4.Range("F7").Value  = 339
$ws4.Range("F9").Value  = 11664
$ws4.Range("F11").Value = 1166
$ws4.Range("F15").Value = 5773
$ws4.Range("F18").Value = 184
$ws4.Range("F19").Value = 18
